# AuditoriaEq3byEq2.xlsx - "Checklist" sheet
# Rows 37-43 (the "Design" artifact section) were re-graded from
# "NC" (Nao Conforme) to "C" (Conforme). All of the dependent totals,
# percentages and COUNTIF summaries further down the sheet recalculate
# automatically from these source cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E37:E43").Value = "C"

# Restore the view/selection state that was active when the workbook was
# last saved (scrolled down to row 37, cell F36 selected).
$win = $excel.ActiveWindow
$win.ScrollRow = 37
$win.ScrollColumn = 1
$ws.Range("F36").Select()

$wb.Save()
